$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '66.893.25'
Set-TextValue 'E2' '  -0.06%  '
Set-TextValue 'D3' '3.074.11'
Set-TextValue 'E3' '  -1.36%  '
Set-TextValue 'E4' '  +0.25%  '
Set-TextValue 'D5' '577.27'
Set-TextValue 'E5' '  -0.15%  '
Set-TextValue 'D6' '168.35'
Set-TextValue 'E6' '  -2.36%  '
Set-TextValue 'E7' '  +0.05%  '
Set-TextValue 'D8' '3.069.38'
Set-TextValue 'E8' '  -1.44%  '
Set-TextValue 'E9' '  -1.99%  '
Set-TextValue 'D10' '6.39'
Set-TextValue 'E10' '  -1.13%  '
Set-TextValue 'E11' '  -2.24%  '
Set-TextValue 'E12' '  -2.78%  '
Set-TextValue 'E13' '  -2.35%  '
Set-TextValue 'D14' '35.85'
Set-TextValue 'E14' '  -3.84%  '
Set-TextValue 'E15' '  -1.60%  '
Set-TextValue 'B16' 'WrappedBTC'
Set-TextValue 'C16' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D16' '66.930.71'
Set-TextValue 'E16' '  +0.17%  '
Set-TextValue 'B17' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C17' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D17' '3.588.66'
Set-TextValue 'E17' '  -1.15%  '
Set-TextValue 'D18' '7.02'
Set-TextValue 'E18' '  -1.82%  '
Set-TextValue 'D19' '16.90'
Set-TextValue 'E19' '  +3.54%  '
Set-TextValue 'D20' '3.091.27'
Set-TextValue 'E20' '  -0.70%  '
Set-TextValue 'D21' '489.23'
Set-TextValue 'E21' '  +2.49%  '
Set-TextValue 'E22' '  -3.57%  '
Set-TextValue 'E23' '  -4.40%  '
Set-TextValue 'D24' '82.84'
Set-TextValue 'E24' '  -1.43%  '
Set-TextValue 'D25' '12.80'
Set-TextValue 'E25' '  -5.04%  '
Set-TextValue 'D26' '2.22'
Set-TextValue 'E26' '  -3.78%  '
Set-TextValue 'D27' '10.31'
Set-TextValue 'E27' '  +2.85%  '
Set-TextValue 'E28' '  +0.01%  '
Set-TextValue 'E29' '  -1.07%  '
Set-TextValue 'E30' '  -5.38%  '
Set-TextValue 'E31' '  -2.08%  '
Set-TextValue 'D32' '27.61'
Set-TextValue 'E32' '  -3.67%  '
Set-TextValue 'E33' '  -2.41%  '
Set-TextValue 'D34' '0.0₃0909'
Set-TextValue 'E34' '  -4.36%  '
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  +0.27%  '
Set-TextValue 'D36' '5.64'
Set-TextValue 'E36' '  -4.04%  '
Set-TextValue 'E37' '  -2.75%  '
Set-TextValue 'D38' '46.92'
Set-TextValue 'E38' '  -0.23%  '
Set-TextValue 'E39' '  +0.67%  '
Set-TextValue 'E40' '  -5.28%  '
Set-TextValue 'E41' '  -2.68%  '
Set-TextValue 'E42' '  -4.08%  '
Set-TextValue 'D43' '2.760.13'
Set-TextValue 'E43' '  -1.98%  '
Set-TextValue 'D44' '373.82'
Set-TextValue 'E44' '  -1.96%  '
Set-TextValue 'B45' 'VeChain'
Set-TextValue 'C45' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D45' '0.0347'
Set-TextValue 'E45' '  -2.56%  '
Set-TextValue 'B46' 'Monero'
Set-TextValue 'C46' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D46' '135.85'
Set-TextValue 'E46' '  -0.17%  '
Set-TextValue 'E47' '  -5.17%  '
Set-TextValue 'E48' '  -0.05%  '
Set-TextValue 'D49' '24.50'
Set-TextValue 'E49' '  -1.92%  '
Set-TextValue 'D50' '2.16'
Set-TextValue 'E50' '  -1.84%  '
Set-TextValue 'E51' '  -1.87%  '
